# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the a021463a-128e-414b-add4-92b0e2388662.md row (row 6) across the three
# report sheets, reflecting a freshly generated handoff xliff.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-26 16:42:40"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-26 16:42:36"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-26 16:42:40"
